$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 205, shifting rows 205:313 down to 206:314.
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new data point.
$ws.Cells.Item(205, 1).Value = 8
$ws.Cells.Item(205, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(205, 3).Value = "Coquimbo"
$ws.Cells.Item(205, 4).Value = 44813
$ws.Cells.Item(205, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(205, 5).Value = 4
$ws.Cells.Item(205, 6).Value = 100112012
$ws.Cells.Item(205, 7).Value = "Espinaca"
$ws.Cells.Item(205, 8).Value = "Sin especificar"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 2400
$ws.Cells.Item(205, 11).Value = 450
$ws.Cells.Item(205, 12).Value = 500
$ws.Cells.Item(205, 13).Value = 475
$ws.Cells.Item(205, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(205, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(205, 16).Value = 950
$ws.Cells.Item(205, 17).Value = 0.5
$ws.Cells.Item(205, 18).Value = "Hortaliza"
